$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("I have diarrhea", "下痢です。|げりです。"),
    @("I am constipated.", "便秘です。|べんぴです。"),
    @("I have my period.", "生理です。|せいりです。"),
    @("I have hay fever.", "花粉症です。|かふんしょうです。"),
    @("I have an allergy to...", "（～に）アレルギーがあります。"),
    @("I have a bad tooth.", "虫歯があります。|むしばがあります。"),
    @("I sneeze.", "くしゃみが出ます。|くしゃみがでます。"),
    @("I have a runny nose.", "鼻水が出ます。|はなみずがでます。"),
    @("My back itches.", "背中がかゆいです。|せなかがかゆいです。"),
    @("I have rashes.", "発疹があります。|はっしんがあります。"),
    @("I feel dizzy.", "めまいがします。"),
    @("I threw up.", "吐きました。|はきました。"),
    @("I am not feeling well.", "気分が悪いです。|きぶんがわるいです。"),
    @("I burned myself.", "やけどをしました。"),
    @("I broke my leg.", "足を骨折しました。|あしをこっせつしました。"),
    @("I hurt myself.", "けがをしました。"),
    @("physician", "内科|ないか"),
    @("dermatologist", "皮膚科|ひふか"),
    @("surgeon", "外科|げか"),
    @("obstetrician and gynecologist", "産婦人科|さんふじんか"),
    @("orthopedic surgeon", "整形外科|せいけいげか"),
    @("ophthalmologist", "眼科|がんか"),
    @("dentist", "歯科|しか"),
    @("otorhinolaryngologist; ENT doctor", "耳鼻科|じびか"),
    @("antibiotic", "抗生物質|こうせいぶっしつ"),
    @("X-ray", "レントゲン"),
    @("operation", "手術|しゅじゅつ"),
    @("injection", "注射|ちゅうしゃ"),
    @("thermometer", "体温計|たいおんけい"),
    @("intravenous feeding", "点滴|てんてき")
)

$startRow = 55
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

